# Slide 4 ("Differential Gene Expression Analysis") - the "TextBox 2" shape
# is rewritten: new wording, an added bolded "Option 2" lead-in on a new
# paragraph, a smaller font size, and the textbox shrinks (spAutoFit) to
# its new, shorter content.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("TextBox 2")

# Replace the body copy with the two new paragraphs. A literal `` `r ``
# (carriage return) inside TextRange.Text starts a new paragraph, matching
# native PowerPoint COM behaviour.
$tr = $sh.TextFrame.TextRange
$tr.Text = "No official statistical method to find DGE across conditions.`rOption 2: Qualitative Visualization"
$tr.Font.Size = 16

# Bold just the "Option 2" run that leads the second paragraph.
$optionRun = $tr.Characters(63, 8)
$optionRun.Font.Bold = $true

# Resize/reposition the textbox so it hugs the now-shorter text (matches
# the author's manual nudge in the authored deck: off y 828675->694577 EMU,
# ext cy 646331->584775 EMU). Shape.Top/Height are expressed in points
# (1 pt = 12700 EMU); the literals below are chosen so the point -> EMU
# round trip lands exactly on the target EMU values.
$sh.Top = 54.69110298220257
$sh.Height = 46.04527559055118
